# DOC125-Controle e Rastreio para Bobina de Cartao Ponto
# "Correcao - Rastreio de Bubina"
#
# 1) Merge the "MARCA: Thega DML" run fragments (and drop the spell-check
#    proofErr bookends around "Thega") into a single plain run.
# 2) Fill in the missing "Data Remocao" date (19/09/2022) for the
#    "Bobina 000003" row of the tracking table.
# 3) Trim the top/bottom table cell margins from the three custom table
#    styles ("a", "a0", "a1") used by the body/header/footer tables.

$d = $word.ActiveDocument

# --- 1) ": Thega DML" --------------------------------------------------
# Collapses ": " + "Thega" (spell-checked) + " DML" into one run with the
# same rPr (sz 20 / szCs 20), removing the proofErr spell-check markers.
$d.Content.Find.Execute(": Thega DML", $false, $false, $false, $false, `
    $false, $true, 1, $false, ": Thega DML", 2) | Out-Null

# --- 2) Missing removal date for "Bobina 000003" -----------------------
$table = $d.Tables.Item(1)
$rowCount = $table.Rows.Count
for ($r = 1; $r -le $rowCount; $r++) {
    $label = $table.Cell($r, 1).Range.Text
    if ($label -like "Bobina 000003*") {
        $cell = $table.Cell($r, 5)
        $rng = $cell.Range
        $rng.InsertBefore("19/09/2022")
        $rng.Font.Color = 0
        $rng.Font.Size = 10
        $rng.Font.SizeBi = 10
        break
    }
}

# --- 3) Drop top/bottom cell margins on table styles a / a0 / a1 -------
# (Styles.Item by name ("a0"/"a1") is unreliable for these anonymous
# custom table styles -- several of them report the same NameLocal, so
# walk the collection and match on Type/BuiltIn/Name instead.)
for ($i = 1; $i -le $d.Styles.Count; $i++) {
    $style = $d.Styles.Item($i)
    if ($style.Type -eq 3 -and $style.BuiltIn -eq $false -and $style.NameLocal -eq "a") {
        $tblFmt = $style.Table
        $tblFmt.TopPadding = 0
        $tblFmt.BottomPadding = 0
    }
}
